# Applies the Dec 14, 2022 GitHub Actions data refresh to the crypto price sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> column letter -> new text value.
$updates = @{
    2 = @{ 'D' = '273.53'; 'F' = '14-12-2022'; 'G' = '0' }
    3 = @{ 'D' = '22.99'; 'F' = '14-12-2022'; 'G' = '0' }
    4 = @{ 'D' = '6.511'; 'F' = '14-12-2022'; 'G' = '0' }
    5 = @{ 'D' = '0.06244'; 'F' = '14-12-2022'; 'G' = '0' }
    6 = @{ 'D' = '3.649'; 'F' = '14-12-2022'; 'G' = '0' }
    7 = @{ 'D' = '6.682'; 'F' = '14-12-2022'; 'G' = '0' }
    8 = @{ 'D' = '1.377'; 'F' = '14-12-2022'; 'G' = '0' }
    9 = @{ 'D' = '0.8333'; 'F' = '14-12-2022'; 'G' = '0' }
    10 = @{ 'D' = '0.01387'; 'F' = '14-12-2022'; 'G' = '0' }
    11 = @{ 'D' = '0.1605'; 'F' = '14-12-2022'; 'G' = '0' }
    12 = @{ 'D' = '0.08285'; 'F' = '14-12-2022'; 'G' = '0' }
    13 = @{ 'D' = '0.03435'; 'F' = '14-12-2022'; 'G' = '0' }
    14 = @{ 'D' = '0.03171'; 'F' = '14-12-2022'; 'G' = '0' }
    15 = @{ 'D' = '0.09339'; 'F' = '14-12-2022'; 'G' = '0' }
    16 = @{ 'D' = '3.840'; 'F' = '14-12-2022'; 'G' = '0' }
    17 = @{ 'D' = '0.001660'; 'F' = '14-12-2022'; 'G' = '0' }
    18 = @{ 'D' = '0.04756'; 'F' = '14-12-2022'; 'G' = '0' }
    19 = @{ 'D' = '0.006289'; 'F' = '14-12-2022'; 'G' = '0' }
    20 = @{ 'D' = '0.005707'; 'F' = '14-12-2022'; 'G' = '0' }
    21 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    22 = @{ 'D' = '0.0001498'; 'F' = '14-12-2022'; 'G' = '0' }
    23 = @{ 'D' = '3.712'; 'F' = '14-12-2022'; 'G' = '0' }
    24 = @{ 'D' = '2.391'; 'F' = '14-12-2022'; 'G' = '0' }
    25 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    26 = @{ 'D' = '0.1253'; 'F' = '14-12-2022'; 'G' = '0' }
    27 = @{ 'D' = '0.0002691'; 'F' = '14-12-2022'; 'G' = '0' }
    28 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    29 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    30 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    31 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    32 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    33 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    34 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    35 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    36 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    37 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    38 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    39 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    40 = @{ 'D' = '0.04714'; 'F' = '14-12-2022'; 'G' = '0' }
    41 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    42 = @{ 'B' = 'BKEXToken'; 'C' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'; 'D' = '0.1165'; 'E' = '41BKEXTokenBKK'; 'F' = '14-12-2022'; 'G' = '0' }
    43 = @{ 'B' = 'CEJI'; 'C' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'; 'D' = '0.003535'; 'E' = '42CEJICEJIWorstin24h'; 'F' = '14-12-2022'; 'G' = '0' }
    44 = @{ 'D' = '0.01175'; 'F' = '14-12-2022'; 'G' = '0' }
    45 = @{ 'D' = '0.00006270'; 'F' = '14-12-2022'; 'G' = '0' }
    46 = @{ 'F' = '14-12-2022'; 'G' = '0' }
    47 = @{ 'D' = '0.00000000749'; 'F' = '14-12-2022'; 'G' = '0' }
    48 = @{ 'D' = '0.7961'; 'F' = '14-12-2022'; 'G' = '0' }
    49 = @{ 'D' = '0.002100'; 'F' = '14-12-2022'; 'G' = '0' }
    50 = @{ 'D' = '0.00001398'; 'E' = '49CryptobidCoinCBC'; 'F' = '14-12-2022'; 'G' = '0' }
    51 = @{ 'D' = '0.01239'; 'F' = '14-12-2022'; 'G' = '0' }
}

foreach ($row in $updates.Keys) {
    $rowUpdates = $updates[$row]
    foreach ($col in $rowUpdates.Keys) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        # Force text storage so values stay exact strings (e.g. "273.53", "14-12-2022")
        # instead of being auto-converted to numbers/dates by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = $rowUpdates[$col]
        # Drop back to the Normal style so no stray quote-prefix/formatting is left behind.
        $cell.Style = "Normal"
    }
}
